# Apply updated probability values to Sheet1 ("team_specific_matrix" style data)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Cells.Item(2, 2).Value = 0.2044198895027624
$ws.Cells.Item(2, 3).Value = 0.5248618784530387
$ws.Cells.Item(2, 10).Value = 0.03038674033149171
$ws.Cells.Item(2, 16).Value = 0.138121546961326
$ws.Cells.Item(2, 19).Value = 0.1022099447513812

# Row 3
$ws.Cells.Item(3, 2).Value = 0.01530612244897959
$ws.Cells.Item(3, 3).Value = 0.01020408163265306
$ws.Cells.Item(3, 10).Value = 0.02551020408163265
$ws.Cells.Item(3, 16).Value = 0.7295918367346939
$ws.Cells.Item(3, 19).Value = 0.2193877551020408

# Row 4
$ws.Cells.Item(4, 10).Value = 0.05357142857142857
$ws.Cells.Item(4, 16).Value = 0.6964285714285714
$ws.Cells.Item(4, 19).Value = 0.25

# Row 6
$ws.Cells.Item(6, 2).Value = 0.05150214592274678
$ws.Cells.Item(6, 4).Value = 0.004291845493562232
$ws.Cells.Item(6, 6).Value = 0.06866952789699571
$ws.Cells.Item(6, 10).Value = 0.3047210300429185
$ws.Cells.Item(6, 15).Value = 0.02145922746781116
$ws.Cells.Item(6, 17).Value = 0.2103004291845494
$ws.Cells.Item(6, 18).Value = 0.06437768240343347
$ws.Cells.Item(6, 19).Value = 0.2746781115879828

# Row 7
$ws.Cells.Item(7, 2).Value = 0.1538461538461539
$ws.Cells.Item(7, 4).Value = 0.03076923076923077
$ws.Cells.Item(7, 6).Value = 0.04102564102564103
$ws.Cells.Item(7, 10).Value = 0.1333333333333333
$ws.Cells.Item(7, 15).Value = 0.01025641025641026
$ws.Cells.Item(7, 17).Value = 0.1692307692307692
$ws.Cells.Item(7, 18).Value = 0.09230769230769231
$ws.Cells.Item(7, 19).Value = 0.3692307692307693

# Row 8
$ws.Cells.Item(8, 2).Value = 0.1347517730496454
$ws.Cells.Item(8, 4).Value = 0.02364066193853428
$ws.Cells.Item(8, 6).Value = 0.05673758865248227
$ws.Cells.Item(8, 10).Value = 0.115839243498818
$ws.Cells.Item(8, 15).Value = 0.007092198581560284
$ws.Cells.Item(8, 17).Value = 0.1607565011820331
$ws.Cells.Item(8, 18).Value = 0.1016548463356974
$ws.Cells.Item(8, 19).Value = 0.3995271867612293

# Row 9
$ws.Cells.Item(9, 2).Value = 0.078125
$ws.Cells.Item(9, 4).Value = 0.01041666666666667
$ws.Cells.Item(9, 6).Value = 0.046875
$ws.Cells.Item(9, 10).Value = 0.1302083333333333
$ws.Cells.Item(9, 15).Value = 0.005208333333333333
$ws.Cells.Item(9, 17).Value = 0.2552083333333333
$ws.Cells.Item(9, 18).Value = 0.0625
$ws.Cells.Item(9, 19).Value = 0.4114583333333333

# Row 10
$ws.Cells.Item(10, 2).Value = 0.1259370314842579
$ws.Cells.Item(10, 4).Value = 0.02848575712143928
$ws.Cells.Item(10, 5).Value = 0.0007496251874062968
$ws.Cells.Item(10, 6).Value = 0.0704647676161919
$ws.Cells.Item(10, 10).Value = 0.1199400299850075
$ws.Cells.Item(10, 15).Value = 0.008245877061469266
$ws.Cells.Item(10, 17).Value = 0.2061469265367316
$ws.Cells.Item(10, 18).Value = 0.07721139430284858
$ws.Cells.Item(10, 19).Value = 0.3628185907046477

# Row 11
$ws.Cells.Item(11, 7).Value = 0.1523809523809524
$ws.Cells.Item(11, 10).Value = 0.09206349206349207
$ws.Cells.Item(11, 11).Value = 0.2253968253968254
$ws.Cells.Item(11, 12).Value = 0.5047619047619047
$ws.Cells.Item(11, 19).Value = 0.0253968253968254

# Row 12
$ws.Cells.Item(12, 7).Value = 0.7619047619047619
$ws.Cells.Item(12, 10).Value = 0.1845238095238095
$ws.Cells.Item(12, 11).Value = 0.005952380952380952
$ws.Cells.Item(12, 12).Value = 0.03571428571428571
$ws.Cells.Item(12, 19).Value = 0.0119047619047619

# Row 13
$ws.Cells.Item(13, 7).Value = 0.7241379310344828
$ws.Cells.Item(13, 10).Value = 0.1724137931034483
$ws.Cells.Item(13, 19).Value = 0.103448275862069

# Row 14
$ws.Cells.Item(14, 7).Value = 0.75
$ws.Cells.Item(14, 19).Value = 0.25

# Row 15
$ws.Cells.Item(15, 6).Value = 0.04739336492890995
$ws.Cells.Item(15, 8).Value = 0.1137440758293839
$ws.Cells.Item(15, 9).Value = 0.08530805687203792
$ws.Cells.Item(15, 10).Value = 0.4265402843601896
$ws.Cells.Item(15, 11).Value = 0.07582938388625593
$ws.Cells.Item(15, 13).Value = 0.01895734597156398
$ws.Cells.Item(15, 15).Value = 0.06635071090047394
$ws.Cells.Item(15, 19).Value = 0.1658767772511848

# Row 16
$ws.Cells.Item(16, 6).Value = 0.01801801801801802
$ws.Cells.Item(16, 8).Value = 0.1846846846846847
$ws.Cells.Item(16, 9).Value = 0.06756756756756757
$ws.Cells.Item(16, 10).Value = 0.4144144144144144
$ws.Cells.Item(16, 11).Value = 0.1261261261261261
$ws.Cells.Item(16, 13).Value = 0.01801801801801802
$ws.Cells.Item(16, 14).Value = 0.009009009009009009
$ws.Cells.Item(16, 15).Value = 0.02702702702702703
$ws.Cells.Item(16, 19).Value = 0.1351351351351351

# Row 17
$ws.Cells.Item(17, 6).Value = 0.02365591397849462
$ws.Cells.Item(17, 8).Value = 0.178494623655914
$ws.Cells.Item(17, 9).Value = 0.07956989247311828
$ws.Cells.Item(17, 10).Value = 0.4559139784946237
$ws.Cells.Item(17, 11).Value = 0.0989247311827957
$ws.Cells.Item(17, 13).Value = 0.01075268817204301
$ws.Cells.Item(17, 14).Value = 0.002150537634408602
$ws.Cells.Item(17, 15).Value = 0.04731182795698925
$ws.Cells.Item(17, 19).Value = 0.1032258064516129

# Row 18
$ws.Cells.Item(18, 6).Value = 0.01058201058201058
$ws.Cells.Item(18, 8).Value = 0.164021164021164
$ws.Cells.Item(18, 9).Value = 0.06878306878306878
$ws.Cells.Item(18, 10).Value = 0.4232804232804233
$ws.Cells.Item(18, 11).Value = 0.07936507936507936
$ws.Cells.Item(18, 13).Value = 0.005291005291005291
$ws.Cells.Item(18, 15).Value = 0.1005291005291005
$ws.Cells.Item(18, 19).Value = 0.1481481481481481

# Row 19
$ws.Cells.Item(19, 6).Value = 0.02427564604541895
$ws.Cells.Item(19, 8).Value = 0.1918559122944401
$ws.Cells.Item(19, 9).Value = 0.08692247454972592
$ws.Cells.Item(19, 10).Value = 0.3625685199686766
$ws.Cells.Item(19, 11).Value = 0.1049334377447142
$ws.Cells.Item(19, 13).Value = 0.01252936570086139
$ws.Cells.Item(19, 14).Value = 0.001566170712607674
$ws.Cells.Item(19, 15).Value = 0.07909162098668755
$ws.Cells.Item(19, 19).Value = 0.1362568519968677
